$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 300.27274
$ws.Range("I4").Value = 144.77777
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 144.77777
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -30.77777
$ws.Range("N4").Value = -1228

$ws.Range("H86").Value = 1074.0769
$ws.Range("I86").Value = 868.5714
$ws.Range("J86").Value = 1313.8334
$ws.Range("K86").Value = 868.5714
$ws.Range("L86").Value = 1313.8334
$ws.Range("M86").Value = 254.4286
$ws.Range("N86").Value = -3559.8334

$ws.Range("H89").Value = 1074.0769
$ws.Range("I89").Value = 868.5714
$ws.Range("J89").Value = 1313.8334
$ws.Range("K89").Value = 4342.857
$ws.Range("L89").Value = 6569.166999999999
$ws.Range("M89").Value = 1273.143
$ws.Range("N89").Value = -17801.167

$ws.Range("H98").Value = 3651.7856
$ws.Range("I98").Value = 1960.25
$ws.Range("J98").Value = 5907.1665
$ws.Range("K98").Value = 1960.25
$ws.Range("L98").Value = 5907.1665
$ws.Range("M98").Value = -462.25
$ws.Range("N98").Value = -8903.166499999999

$ws.Range("H112").Value = 1381.0103
$ws.Range("I112").Value = 450
$ws.Range("J112").Value = 1400.6105
$ws.Range("K112").Value = 1350
$ws.Range("L112").Value = 4201.8315
$ws.Range("M112").Value = -242
$ws.Range("N112").Value = -6417.8315

$ws.Range("H122").Value = 3651.7856
$ws.Range("I122").Value = 1960.25
$ws.Range("J122").Value = 5907.1665
$ws.Range("K122").Value = 5880.75
$ws.Range("L122").Value = 17721.4995
$ws.Range("M122").Value = -3430.75
$ws.Range("N122").Value = -22621.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = 16

$ws.Range("H5").Value = 280
$ws.Range("I5").Value = 280
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 280
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -168

$ws.Range("H74").Value = 1043.7858
$ws.Range("I74").Value = 884.5
$ws.Range("K74").Value = 884.5
$ws.Range("M74").Value = -10.5

$ws.Range("H77").Value = 1043.7858
$ws.Range("I77").Value = 884.5
$ws.Range("K77").Value = 4422.5
$ws.Range("M77").Value = -54.5

$ws.Range("H132").Value = 2762.25
$ws.Range("I132").Value = 2233.2173
$ws.Range("J132").Value = 4114.222
$ws.Range("K132").Value = 6699.651899999999
$ws.Range("L132").Value = 12342.666
$ws.Range("M132").Value = -4169.651899999999
$ws.Range("N132").Value = -17402.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 280
$ws.Range("I4").Value = 280
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 280
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -165

$ws.Range("H107").Value = 3862.3333
$ws.Range("I107").Value = 3195.4546
$ws.Range("J107").Value = 4595.9
$ws.Range("K107").Value = 3195.4546
$ws.Range("L107").Value = 4595.9
$ws.Range("M107").Value = -1275.4546
$ws.Range("N107").Value = -8435.9

$ws.Range("H133").Value = 39990
$ws.Range("J133").Value = 39990
$ws.Range("L133").Value = 39990
$ws.Range("N133").Value = -50110

$ws.Range("H134").Value = 4290.1665
$ws.Range("I134").Value = 2787.4285
$ws.Range("J134").Value = 6394
$ws.Range("K134").Value = 8362.2855
$ws.Range("L134").Value = 19182
$ws.Range("M134").Value = -5827.2855
$ws.Range("N134").Value = -24252

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9117
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 9654
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 9654
$ws.Range("M4").Value = -4888
$ws.Range("N4").Value = -9878

$ws.Range("H7").Value = 71.375
$ws.Range("I7").Value = 57.75
$ws.Range("J7").Value = 85
$ws.Range("K7").Value = 57.75
$ws.Range("L7").Value = 85
$ws.Range("M7").Value = 55.25
$ws.Range("N7").Value = -311

$ws.Range("H31").Value = 1756995.9
$ws.Range("J31").Value = 4540.0527
$ws.Range("L31").Value = 4540.0527
$ws.Range("N31").Value = -5130.0527

$ws.Range("H34").Value = 1756995.9
$ws.Range("J34").Value = 4540.0527
$ws.Range("L34").Value = 4540.0527
$ws.Range("N34").Value = -4944.0527

$ws.Range("H74").Value = 16435.166
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 17474.727
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 17474.727
$ws.Range("M74").Value = -4126
$ws.Range("N74").Value = -19222.727

$ws.Range("H77").Value = 16435.166
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 17474.727
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 52424.181
$ws.Range("M77").Value = -10632
$ws.Range("N77").Value = -61160.181

$ws.Range("H132").Value = 1675.196
$ws.Range("I132").Value = 1199.0233
$ws.Range("J132").Value = 4234.625
$ws.Range("K132").Value = 3597.0699
$ws.Range("L132").Value = 12703.875
$ws.Range("M132").Value = -1067.0699
$ws.Range("N132").Value = -17763.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 805.3333
$ws.Range("I4").Value = 76.666664
$ws.Range("J4").Value = 1534
$ws.Range("K4").Value = 229.999992
$ws.Range("L4").Value = 4602
$ws.Range("M4").Value = -117.999992
$ws.Range("N4").Value = -4826

$ws.Range("H18").Value = 650.7692
$ws.Range("I18").Value = 346
$ws.Range("J18").Value = 1666.6666
$ws.Range("K18").Value = 1038
$ws.Range("L18").Value = 4999.9998
$ws.Range("M18").Value = -869
$ws.Range("N18").Value = -5337.9998

$ws.Range("H68").Value = 8225.75
$ws.Range("J68").Value = 10667.667
$ws.Range("L68").Value = 32003.001
$ws.Range("N68").Value = -33625.001

$ws.Range("H71").Value = 8225.75
$ws.Range("J71").Value = 10667.667
$ws.Range("L71").Value = 96009.003
$ws.Range("N71").Value = -104121.003

$ws.Range("H80").Value = 3475.25
$ws.Range("I80").Value = 1499
$ws.Range("J80").Value = 3654.9092
$ws.Range("K80").Value = 4497
$ws.Range("L80").Value = 10964.7276
$ws.Range("M80").Value = -3561
$ws.Range("N80").Value = -12836.7276

$ws.Range("H83").Value = 3475.25
$ws.Range("I83").Value = 1499
$ws.Range("J83").Value = 3654.9092
$ws.Range("K83").Value = 13491
$ws.Range("L83").Value = 32894.1828
$ws.Range("M83").Value = -8811
$ws.Range("N83").Value = -42254.1828

$ws.Range("H106").Value = 3557.1428
$ws.Range("J106").Value = 3557.1428
$ws.Range("L106").Value = 10671.4284
$ws.Range("N106").Value = -12563.4284

$ws.Range("H127").Value = 1780.75
$ws.Range("J127").Value = 1780.75
$ws.Range("L127").Value = 5342.25
$ws.Range("N127").Value = -15262.25

$ws.Range("H131").Value = 1411.125
$ws.Range("I131").Value = 1901.2
$ws.Range("J131").Value = 1188.3636
$ws.Range("K131").Value = 5703.6
$ws.Range("L131").Value = 3565.0908
$ws.Range("M131").Value = -663.6000000000004
$ws.Range("N131").Value = -13645.0908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5013.5386
$ws.Range("I70").Value = 5199.636
$ws.Range("J70").Value = 3990
$ws.Range("K70").Value = 5199.636
$ws.Range("L70").Value = 3990
$ws.Range("M70").Value = -4929.636
$ws.Range("N70").Value = -4530

$ws.Range("H73").Value = 5013.5386
$ws.Range("I73").Value = 5199.636
$ws.Range("J73").Value = 3990
$ws.Range("K73").Value = 5199.636
$ws.Range("L73").Value = 3990
$ws.Range("M73").Value = -4263.636
$ws.Range("N73").Value = -5862

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 221183.56
$ws.Range("J2").Value = 24773.38
$ws.Range("L2").Value = 24773.38
$ws.Range("N2").Value = -24997.38

$ws.Range("H7").Value = 1909.9166
$ws.Range("I7").Value = 1102.1111
$ws.Range("K7").Value = 1102.1111
$ws.Range("M7").Value = -990.1111000000001

$ws.Range("H68").Value = 2366.7917
$ws.Range("I68").Value = 1052.6316
$ws.Range("J68").Value = 7360.6
$ws.Range("K68").Value = 1052.6316
$ws.Range("L68").Value = 7360.6
$ws.Range("M68").Value = -303.6315999999999
$ws.Range("N68").Value = -8858.6

$ws.Range("H71").Value = 2366.7917
$ws.Range("I71").Value = 1052.6316
$ws.Range("J71").Value = 7360.6
$ws.Range("K71").Value = 5263.157999999999
$ws.Range("L71").Value = 36803
$ws.Range("M71").Value = -1519.157999999999
$ws.Range("N71").Value = -44291

$ws.Range("H122").Value = 3142.8572
$ws.Range("I122").Value = 2750
$ws.Range("J122").Value = 4125
$ws.Range("K122").Value = 8250
$ws.Range("L122").Value = 12375
$ws.Range("M122").Value = -5800
$ws.Range("N122").Value = -17275

$ws.Range("H126").Value = 1909.9166
$ws.Range("I126").Value = 1102.1111
$ws.Range("K126").Value = 3306.3333
$ws.Range("M126").Value = -836.3333000000002

$ws.Range("H135").Value = 29764.047
$ws.Range("J135").Value = 29764.047
$ws.Range("L135").Value = 29764.047
$ws.Range("N135").Value = -39904.047

$ws.Range("H136").Value = 2383477.2
$ws.Range("I136").Value = 2779167.8
$ws.Range("J136").Value = 9333.333000000001
$ws.Range("K136").Value = 8337503.399999999
$ws.Range("L136").Value = 27999.999
$ws.Range("M136").Value = -8334953.399999999
$ws.Range("N136").Value = -33099.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1088.7858
$ws.Range("I113").Value = 498.17648
$ws.Range("J113").Value = 2001.5454
$ws.Range("K113").Value = 1494.52944
$ws.Range("L113").Value = 6004.6362
$ws.Range("M113").Value = 675.47056
$ws.Range("N113").Value = -10344.6362

$ws.Range("H122").Value = 528593.7
$ws.Range("I122").Value = 626630.0600000001
$ws.Range("K122").Value = 1879890.18
$ws.Range("M122").Value = -1877440.18

$ws.Range("H126").Value = 3449650.2
$ws.Range("I126").Value = 763.38464
$ws.Range("K126").Value = 2290.15392
$ws.Range("M126").Value = 179.8460800000003

$ws.Range("H140").Value = 62085.6
$ws.Range("J140").Value = 62085.6
$ws.Range("L140").Value = 62085.6
$ws.Range("N140").Value = -72445.60000000001
